$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 1699
$ws.Cells.Item(32, 10).Value = 1330
$ws.Cells.Item(32, 12).Value = 1330
$ws.Cells.Item(32, 14).Value = -1982
$ws.Cells.Item(38, 8).Value = 42.25
$ws.Cells.Item(38, 9).Value = 42.25
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 11).Value = 126.75
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 13).Value = 245.25
$ws.Cells.Item(38, 14).ClearContents()
$ws.Cells.Item(40, 8).Value = 3824
$ws.Cells.Item(40, 10).Value = 5497
$ws.Cells.Item(40, 12).Value = 5497
$ws.Cells.Item(40, 14).Value = -5847
$ws.Cells.Item(41, 8).Value = 797.44446
$ws.Cells.Item(41, 9).Value = 675.8
$ws.Cells.Item(41, 10).Value = 949.5
$ws.Cells.Item(41, 11).Value = 675.8
$ws.Cells.Item(41, 12).Value = 949.5
$ws.Cells.Item(41, 13).Value = -235.8
$ws.Cells.Item(41, 14).Value = -1829.5
$ws.Cells.Item(70, 8).Value = 4142.357
$ws.Cells.Item(70, 10).Value = 4286.9165
$ws.Cells.Item(70, 12).Value = 12860.7495
$ws.Cells.Item(70, 14).Value = -13400.7495
$ws.Cells.Item(73, 8).Value = 4142.357
$ws.Cells.Item(73, 10).Value = 4286.9165
$ws.Cells.Item(73, 12).Value = 12860.7495
$ws.Cells.Item(73, 14).Value = -14732.7495
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 14).ClearContents()
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 12).Value = 0
$ws.Cells.Item(90, 14).ClearContents()
$ws.Cells.Item(106, 8).Value = 2356
$ws.Cells.Item(106, 9).Value = 1867.1578
$ws.Cells.Item(106, 11).Value = 1867.1578
$ws.Cells.Item(106, 13).Value = -1236.1578
$ws.Cells.Item(132, 8).Value = 7775.6772
$ws.Cells.Item(132, 9).Value = 7904.6665
$ws.Cells.Item(132, 10).Value = 3906
$ws.Cells.Item(132, 11).Value = 23713.9995
$ws.Cells.Item(132, 12).Value = 11718
$ws.Cells.Item(132, 13).Value = -21183.9995
$ws.Cells.Item(132, 14).Value = -16778
$ws.Cells.Item(138, 8).Value = 7945.609
$ws.Cells.Item(138, 9).Value = 15000
$ws.Cells.Item(138, 10).Value = 6887.45
$ws.Cells.Item(138, 11).Value = 45000
$ws.Cells.Item(138, 12).Value = 20662.35
$ws.Cells.Item(138, 13).Value = -39860
$ws.Cells.Item(138, 14).Value = -30942.35
$ws.Cells.Item(141, 8).Value = 10834.104
$ws.Cells.Item(141, 9).Value = 2898.625
$ws.Cells.Item(141, 11).Value = 8695.875
$ws.Cells.Item(141, 13).Value = -3515.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1068.76
$ws.Cells.Item(2, 9).Value = 960.0909
$ws.Cells.Item(2, 11).Value = 960.0909
$ws.Cells.Item(2, 13).Value = -847.0909
$ws.Cells.Item(32, 8).Value = 914300.7
$ws.Cells.Item(32, 9).Value = 970278.4
$ws.Cells.Item(32, 11).Value = 970278.4
$ws.Cells.Item(32, 13).Value = -969991.4
$ws.Cells.Item(61, 8).Value = 3848030.2
$ws.Cells.Item(61, 9).Value = 2017.3478
$ws.Cells.Item(61, 10).Value = 33334128
$ws.Cells.Item(61, 11).Value = 2017.3478
$ws.Cells.Item(61, 12).Value = 33334128
$ws.Cells.Item(61, 13).Value = -1805.3478
$ws.Cells.Item(61, 14).Value = -33334552
$ws.Cells.Item(63, 8).Value = 2354.889
$ws.Cells.Item(63, 10).Value = 2457
$ws.Cells.Item(63, 12).Value = 2457
$ws.Cells.Item(63, 14).Value = -3829
$ws.Cells.Item(66, 8).Value = 2354.889
$ws.Cells.Item(66, 10).Value = 2457
$ws.Cells.Item(66, 12).Value = 12285
$ws.Cells.Item(66, 14).Value = -19149
$ws.Cells.Item(74, 8).Value = 995303.2
$ws.Cells.Item(74, 9).Value = 1148227.9
$ws.Cells.Item(74, 11).Value = 1148227.9
$ws.Cells.Item(74, 13).Value = -1147353.9
$ws.Cells.Item(77, 8).Value = 995303.2
$ws.Cells.Item(77, 9).Value = 1148227.9
$ws.Cells.Item(77, 11).Value = 5741139.5
$ws.Cells.Item(77, 13).Value = -5736771.5
$ws.Cells.Item(102, 8).Value = 3613.238
$ws.Cells.Item(102, 9).Value = 3619.9473
$ws.Cells.Item(102, 11).Value = 3619.9473
$ws.Cells.Item(102, 13).Value = -1997.9473
$ws.Cells.Item(116, 8).Value = 1068.76
$ws.Cells.Item(116, 9).Value = 960.0909
$ws.Cells.Item(116, 11).Value = 960.0909
$ws.Cells.Item(116, 13).Value = 1333.9091
$ws.Cells.Item(132, 8).Value = 7384.5
$ws.Cells.Item(132, 9).Value = 3995
$ws.Cells.Item(132, 11).Value = 11985
$ws.Cells.Item(132, 13).Value = -9455
$ws.Cells.Item(136, 8).Value = 3848030.2
$ws.Cells.Item(136, 9).Value = 2017.3478
$ws.Cells.Item(136, 10).Value = 33334128
$ws.Cells.Item(136, 11).Value = 6052.0434
$ws.Cells.Item(136, 12).Value = 100002384
$ws.Cells.Item(136, 13).Value = -3502.0434
$ws.Cells.Item(136, 14).Value = -100007484

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1068.76
$ws.Cells.Item(3, 9).Value = 960.0909
$ws.Cells.Item(3, 11).Value = 960.0909
$ws.Cells.Item(3, 13).Value = -846.0909
$ws.Cells.Item(94, 8).Value = 3945.2
$ws.Cells.Item(94, 9).Value = 3272.5
$ws.Cells.Item(94, 11).Value = 3272.5
$ws.Cells.Item(94, 13).Value = -2821.5
$ws.Cells.Item(128, 8).Value = 6000
$ws.Cells.Item(128, 9).Value = 6000
$ws.Cells.Item(128, 11).Value = 18000
$ws.Cells.Item(128, 13).Value = -15510
$ws.Cells.Item(134, 8).Value = 11909834
$ws.Cells.Item(134, 9).Value = 3100
$ws.Cells.Item(134, 11).Value = 9300
$ws.Cells.Item(134, 13).Value = -6765

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 31322.389
$ws.Cells.Item(16, 9).Value = 38836.07
$ws.Cells.Item(16, 11).Value = 38836.07
$ws.Cells.Item(16, 13).Value = -38549.07
$ws.Cells.Item(31, 8).Value = 1044275.25
$ws.Cells.Item(31, 9).Value = 1044275.25
$ws.Cells.Item(31, 11).Value = 1044275.25
$ws.Cells.Item(31, 13).Value = -1043980.25
$ws.Cells.Item(34, 8).Value = 1044275.25
$ws.Cells.Item(34, 9).Value = 1044275.25
$ws.Cells.Item(34, 11).Value = 1044275.25
$ws.Cells.Item(34, 13).Value = -1044073.25
$ws.Cells.Item(58, 8).Value = 8056670.5
$ws.Cells.Item(58, 9).Value = 9806716
$ws.Cells.Item(58, 11).Value = 9806716
$ws.Cells.Item(58, 13).Value = -9806513
$ws.Cells.Item(86, 8).Value = 9418.040000000001
$ws.Cells.Item(86, 9).Value = 7050
$ws.Cells.Item(86, 11).Value = 7050
$ws.Cells.Item(86, 13).Value = -5927
$ws.Cells.Item(89, 8).Value = 9418.040000000001
$ws.Cells.Item(89, 9).Value = 7050
$ws.Cells.Item(89, 11).Value = 35250
$ws.Cells.Item(89, 13).Value = -29634
$ws.Cells.Item(113, 8).Value = 31322.389
$ws.Cells.Item(113, 9).Value = 38836.07
$ws.Cells.Item(113, 11).Value = 38836.07
$ws.Cells.Item(113, 13).Value = -36666.07
$ws.Cells.Item(132, 8).Value = 2051.9285
$ws.Cells.Item(132, 9).Value = 1906.5897
$ws.Cells.Item(132, 11).Value = 5719.7691
$ws.Cells.Item(132, 13).Value = -3189.7691
$ws.Cells.Item(134, 8).Value = 4849.3335
$ws.Cells.Item(134, 9).Value = 3321.647
$ws.Cells.Item(134, 11).Value = 9964.940999999999
$ws.Cells.Item(134, 13).Value = -7429.940999999999
$ws.Cells.Item(136, 8).Value = 8056670.5
$ws.Cells.Item(136, 9).Value = 9806716
$ws.Cells.Item(136, 11).Value = 29420148
$ws.Cells.Item(136, 13).Value = -29417598

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 5586974.5
$ws.Cells.Item(4, 9).Value = 10125314
$ws.Cells.Item(4, 11).Value = 30375942
$ws.Cells.Item(4, 13).Value = -30375830
$ws.Cells.Item(5, 8).Value = 2949206.2
$ws.Cells.Item(5, 9).Value = 2101400.8
$ws.Cells.Item(5, 10).Value = 7753437.5
$ws.Cells.Item(5, 11).Value = 6304202.399999999
$ws.Cells.Item(5, 12).Value = 23260312.5
$ws.Cells.Item(5, 13).Value = -6304090.399999999
$ws.Cells.Item(5, 14).Value = -23260536.5
$ws.Cells.Item(131, 8).Value = 6065.7
$ws.Cells.Item(131, 10).Value = 9277.333000000001
$ws.Cells.Item(131, 12).Value = 27831.999
$ws.Cells.Item(131, 14).Value = -37911.999
$ws.Cells.Item(135, 8).Value = 2949206.2
$ws.Cells.Item(135, 9).Value = 2101400.8
$ws.Cells.Item(135, 10).Value = 7753437.5
$ws.Cells.Item(135, 11).Value = 18912607.2
$ws.Cells.Item(135, 12).Value = 69780937.5
$ws.Cells.Item(135, 13).Value = -18910072.2
$ws.Cells.Item(135, 14).Value = -69786007.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 35590.1
$ws.Cells.Item(122, 9).Value = 37975.105
$ws.Cells.Item(122, 11).Value = 113925.315
$ws.Cells.Item(122, 13).Value = -111475.315
$ws.Cells.Item(123, 8).Value = 87162.5
$ws.Cells.Item(123, 10).Value = 87162.5
$ws.Cells.Item(123, 12).Value = 87162.5
$ws.Cells.Item(123, 14).Value = -92062.5
$ws.Cells.Item(132, 8).Value = 29630.545
$ws.Cells.Item(132, 9).Value = 16705.285
$ws.Cells.Item(132, 11).Value = 50115.855
$ws.Cells.Item(132, 13).Value = -47585.855

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 1044246
$ws.Cells.Item(132, 9).Value = 2085257.6
$ws.Cells.Item(132, 11).Value = 6255772.800000001
$ws.Cells.Item(132, 13).Value = -6253242.800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 2101.739
$ws.Cells.Item(107, 9).Value = 642.3684
$ws.Cells.Item(107, 10).Value = 9033.75
$ws.Cells.Item(107, 11).Value = 1927.1052
$ws.Cells.Item(107, 12).Value = 27101.25
$ws.Cells.Item(107, 13).Value = -7.105199999999968
$ws.Cells.Item(107, 14).Value = -30941.25
$ws.Cells.Item(126, 8).Value = 3589.9092
$ws.Cells.Item(126, 10).Value = 3333.3333
$ws.Cells.Item(126, 12).Value = 9999.999899999999
$ws.Cells.Item(126, 14).Value = -14939.9999
$ws.Cells.Item(132, 8).Value = 3790384.8
$ws.Cells.Item(132, 9).Value = 4067427.5
$ws.Cells.Item(132, 10).Value = 4133
$ws.Cells.Item(132, 11).Value = 12202282.5
$ws.Cells.Item(132, 12).Value = 12399
$ws.Cells.Item(132, 13).Value = -12199752.5
$ws.Cells.Item(132, 14).Value = -17459
$ws.Cells.Item(136, 8).Value = 3148246.8
$ws.Cells.Item(136, 9).Value = 1451130.8
$ws.Cells.Item(136, 11).Value = 4353392.4
$ws.Cells.Item(136, 13).Value = -4350842.4

"Applied 226 cell updates across 8 sheets"